$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header D1: "System Type" -> "Asset Name"
$ws.Range("D1").Value = "Asset Name"

# Update row 2 values
$ws.Range("A2").Value = "Capgemini"
$ws.Range("B2").Value = 45690.397916666669
$ws.Range("C2").Value = 45693.791666666664
$ws.Range("D2").Value = "Warehouse Management System"

# Column widths (engine stores width quantized to 1/6 pt on top of a
# +5/6 padding baked into the ColumnWidth->xml-width conversion, so we
# back out the padding to land as close as possible on the target xml width)
$ws.Columns.Item(2).ColumnWidth = 19.072916666666668
$ws.Columns.Item(3).ColumnWidth = 22.436197916666668
$ws.Columns.Item(4).ColumnWidth = 35.166666666666664

# Update selection
$ws.Range("D12").Select()
